$d = $word.ActiveDocument

# ---- Part 1: "table X" caption text -> "table X 1" (Thai: "ตาราง … " -> "ตารางที่ 1 ") ----
$para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Activity Diagram*") { $para = $p; break }
}
$full = $d.Range($para.Range.Start, $para.Range.End - 1)
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r w:rsidRPr="00C60201"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:hint="cs"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:cs/></w:rPr><w:t>ตารา</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:hint="cs"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:cs/></w:rPr><w:t xml:space="preserve">งที่ </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="003634DB"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Activity Diagram</w:t></w:r></w:p>'
[void]$full.InsertXML($xml1)

# ---- Part 2: drop w:hint="cs" from the paragraph-mark rPr of the "2.7.1" version row ----
$t = $d.Tables.Item(1)
$row = $t.Rows.Item(7)

$c3 = $row.Cells.Item(3)
$p3 = $c3.Range.Paragraphs.Item(1)
$rng3 = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$xml3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="72239EE8" w14:textId="09D31012" w:rsidR="00525A6F" w:rsidRDefault="00525A6F" w:rsidP="00525A6F"><w:pPr><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:cs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:cs/></w:rPr><w:t>แก้ไข</w:t></w:r></w:p>'
[void]$rng3.InsertXML($xml3)

$c4 = $row.Cells.Item(4)
$p4 = $c4.Range.Paragraphs.Item(1)
$rng4 = $d.Range($p4.Range.Start, $p4.Range.End - 1)
$xml4 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="5B942131" w14:textId="011C4202" w:rsidR="00525A6F" w:rsidRDefault="00525A6F" w:rsidP="00525A6F"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:cs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:cs/></w:rPr><w:t xml:space="preserve">วริศรา </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>(D)</w:t></w:r></w:p>'
[void]$rng4.InsertXML($xml4)

$c5 = $row.Cells.Item(5)
$p5 = $c5.Range.Paragraphs.Item(1)
$rng5 = $d.Range($p5.Range.Start, $p5.Range.End - 1)
$xml5 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="5CCA38D4" w14:textId="341CF46F" w:rsidR="00525A6F" w:rsidRDefault="00525A6F" w:rsidP="00525A6F"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:cs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:cs/></w:rPr><w:t xml:space="preserve">วิรัตน์ </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>(TL)</w:t></w:r></w:p>'
[void]$rng5.InsertXML($xml5)
